$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (row 2, 4, 7 test data refreshed for the "26 june" run)
$ws.Range("A2").Value = "mtest0626@gmail.com"
$ws.Range("B2").Value = "Mtest@0623"

$ws.Range("A4").Value = "mtest0626b@gmail.com"
$ws.Range("B4").Value = "Mtest@0623"
$ws.Range("C4").Value = "Valid"

$ws.Range("A7").Value = "mtest0626c@gmail.com"
$ws.Range("B7").Value = "Mtest@0626"

# Turn the refreshed username/password cells into hyperlinks
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:mtest0626@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:mtest0626@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:mtest0626b@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:mtest0626b@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:mtest0626c@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:mtest0626c@gmail.com")

# Move the active selection as left by the author at save time
$ws.Range("C12").Select()
